# Auto-generated edit script applying numeric corrections to H:N price/profit
# columns across multiple sheets, per the scheduled runner update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 21678.574
$ws.Range("I132").Value = 3049.861
$ws.Range("J132").Value = 82645.27
$ws.Range("K132").Value = 9149.582999999999
$ws.Range("L132").Value = 247935.81
$ws.Range("M132").Value = -6619.582999999999
$ws.Range("N132").Value = -252995.81

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1985.05
$ws.Range("I74").Value = 1563.303
$ws.Range("K74").Value = 1563.303
$ws.Range("M74").Value = -689.3030000000001
$ws.Range("H77").Value = 1985.05
$ws.Range("I77").Value = 1563.303
$ws.Range("K77").Value = 7816.515
$ws.Range("M77").Value = -3448.515
$ws.Range("H102").Value = 40286.43
$ws.Range("I102").Value = 17112.223
$ws.Range("J102").Value = 82000
$ws.Range("K102").Value = 17112.223
$ws.Range("L102").Value = 82000
$ws.Range("M102").Value = -15490.223
$ws.Range("N102").Value = -85244
$ws.Range("H122").Value = 1975.5
$ws.Range("I122").Value = 1608.5454
$ws.Range("K122").Value = 4825.6362
$ws.Range("M122").Value = -2375.6362

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 43333.332
$ws.Range("J52").Value = 43333.332
$ws.Range("L52").Value = 43333.332
$ws.Range("N52").Value = -43921.332
$ws.Range("H99").Value = 2124.7144
$ws.Range("I99").Value = 2269.5557
$ws.Range("J99").Value = 2016.0834
$ws.Range("K99").Value = 2269.5557
$ws.Range("L99").Value = 2016.0834
$ws.Range("M99").Value = -771.5556999999999
$ws.Range("N99").Value = -5012.0834
$ws.Range("H126").Value = 2124.7144
$ws.Range("I126").Value = 2269.5557
$ws.Range("J126").Value = 2016.0834
$ws.Range("K126").Value = 6808.6671
$ws.Range("L126").Value = 6048.2502
$ws.Range("M126").Value = -4338.6671
$ws.Range("N126").Value = -10988.2502
$ws.Range("H132").Value = 541591.4
$ws.Range("I132").Value = 1756.8823
$ws.Range("J132").Value = 1561278.8
$ws.Range("K132").Value = 5270.6469
$ws.Range("L132").Value = 4683836.4
$ws.Range("M132").Value = -2740.6469
$ws.Range("N132").Value = -4688896.4
$ws.Range("H134").Value = 1222716.9
$ws.Range("I134").Value = 742720.5600000001
$ws.Range("J134").Value = 3502699.5
$ws.Range("K134").Value = 2228161.68
$ws.Range("L134").Value = 10508098.5
$ws.Range("M134").Value = -2225626.68
$ws.Range("N134").Value = -10513168.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 25851
$ws.Range("I5").Value = 25851
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 77553
$ws.Range("L5").Value = 0
$ws.Range("N5").Value = -77441
$ws.Range("M5").ClearContents()
$ws.Range("H18").Value = 38166.25
$ws.Range("I18").Value = 43575.715
$ws.Range("K18").Value = 130727.145
$ws.Range("M18").Value = -130558.145
$ws.Range("H99").Value = 8064
$ws.Range("I99").Value = 8064
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 24192
$ws.Range("L99").Value = 0
$ws.Range("N99").Value = -21946
$ws.Range("M99").ClearContents()
$ws.Range("H100").Value = 4380
$ws.Range("J100").Value = 4644.4443
$ws.Range("L100").Value = 13933.3329
$ws.Range("N100").Value = -15555.3329
$ws.Range("H103").Value = 2514.2856
$ws.Range("I103").Value = 2150
$ws.Range("K103").Value = 6450
$ws.Range("M103").Value = -5571
$ws.Range("H106").Value = 4650
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 4650
$ws.Range("K106").Value = 0
$ws.Range("M106").Value = 13950
$ws.Range("N106").Value = -15842
$ws.Range("L106").ClearContents()
$ws.Range("H109").Value = 3038.9167
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 3038.9167
$ws.Range("K109").Value = 0
$ws.Range("M109").Value = 9116.750100000001
$ws.Range("N109").Value = -11196.7501
$ws.Range("L109").ClearContents()
$ws.Range("H112").Value = 13336979
$ws.Range("J112").Value = 3906.4285
$ws.Range("L112").Value = 11719.2855
$ws.Range("N112").Value = -13935.2855
$ws.Range("H115").Value = 4120.1904
$ws.Range("I115").Value = 864
$ws.Range("J115").Value = 4462.9473
$ws.Range("K115").Value = 2592
$ws.Range("L115").Value = 13388.8419
$ws.Range("M115").Value = -1417
$ws.Range("N115").Value = -15738.8419
$ws.Range("H118").Value = 3771.1428
$ws.Range("I118").Value = 1250
$ws.Range("J118").Value = 3897.2
$ws.Range("K118").Value = 3750
$ws.Range("L118").Value = 11691.6
$ws.Range("M118").Value = -2507
$ws.Range("N118").Value = -14177.6
$ws.Range("H121").Value = 344793.56
$ws.Range("I121").Value = 327.5
$ws.Range("J121").Value = 620366.4
$ws.Range("K121").Value = 982.5
$ws.Range("L121").Value = 1861099.2
$ws.Range("M121").Value = 327.5
$ws.Range("N121").Value = -1863719.2
$ws.Range("H122").Value = 2694.9805
$ws.Range("I122").Value = 575.12823
$ws.Range("K122").Value = 5176.154070000001
$ws.Range("M122").Value = -2726.154070000001
$ws.Range("H125").Value = 1505044.5
$ws.Range("I125").Value = 10002000
$ws.Range("J125").Value = 5581.7646
$ws.Range("K125").Value = 30006000
$ws.Range("L125").Value = 16745.2938
$ws.Range("M125").Value = -30001080
$ws.Range("N125").Value = -26585.2938
$ws.Range("H129").Value = 137831.4
$ws.Range("J129").Value = 1950.3846
$ws.Range("L129").Value = 5851.1538
$ws.Range("N129").Value = -15851.1538
$ws.Range("H130").Value = 51423.168
$ws.Range("J130").Value = 2002.25
$ws.Range("L130").Value = 6006.75
$ws.Range("N130").Value = -16046.75
$ws.Range("H131").Value = 7208.1113
$ws.Range("I131").Value = 10594.7
$ws.Range("J131").Value = 2974.875
$ws.Range("K131").Value = 31784.1
$ws.Range("L131").Value = 8924.625
$ws.Range("M131").Value = -26744.1
$ws.Range("N131").Value = -19004.625
$ws.Range("H132").Value = 2852.95
$ws.Range("I132").Value = 1157.6364
$ws.Range("J132").Value = 4925
$ws.Range("K132").Value = 10418.7276
$ws.Range("L132").Value = 44325
$ws.Range("M132").Value = -7888.7276
$ws.Range("N132").Value = -49385
$ws.Range("H134").Value = 50207276
$ws.Range("I134").Value = 50207276
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 150621828
$ws.Range("L134").Value = 0
$ws.Range("N134").Value = -150616758
$ws.Range("M134").ClearContents()
$ws.Range("H135").Value = 25851
$ws.Range("I135").Value = 25851
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 232659
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = -230124
$ws.Range("M135").ClearContents()
$ws.Range("H136").Value = 26317690
$ws.Range("I136").Value = 33335006
$ws.Range("K136").Value = 100005018
$ws.Range("M136").Value = -99999918
$ws.Range("H137").Value = 17742.615
$ws.Range("I137").Value = 3840
$ws.Range("J137").Value = 39986.8
$ws.Range("K137").Value = 11520
$ws.Range("L137").Value = 119960.4
$ws.Range("M137").Value = -6420
$ws.Range("N137").Value = -130160.4
$ws.Range("H138").Value = 2912.8572
$ws.Range("I138").Value = 2912.8572
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 8738.571599999999
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = -3598.571599999999
$ws.Range("M138").ClearContents()
$ws.Range("H139").Value = 223626.14
$ws.Range("I139").Value = 223626.14
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 670878.42
$ws.Range("L139").Value = 0
$ws.Range("N139").Value = -665738.42
$ws.Range("M139").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 55561724
$ws.Range("I132").Value = 142866260
$ws.Range("J132").Value = 4293.364
$ws.Range("K132").Value = 428598780
$ws.Range("L132").Value = 12880.092
$ws.Range("M132").Value = -428596250
$ws.Range("N132").Value = -17940.092

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4763.0713
$ws.Range("I40").Value = 3714.4285
$ws.Range("J40").Value = 5811.7144
$ws.Range("K40").Value = 3714.4285
$ws.Range("L40").Value = 5811.7144
$ws.Range("M40").Value = -3578.4285
$ws.Range("N40").Value = -6083.7144

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1361498.8
$ws.Range("I122").Value = 2198718
$ws.Range("J122").Value = 1017.5
$ws.Range("K122").Value = 6596154
$ws.Range("L122").Value = 2198718
$ws.Range("M122").Value = -6593704
$ws.Range("N122").Value = -7952.5
$ws.Range("H136").Value = 403157.16
$ws.Range("I136").Value = 583918.2
$ws.Range("K136").Value = 1751754.6
$ws.Range("M136").Value = -1749204.6

